$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row at row 4 for "Showroom Paused" (shifts existing rows 4+ down by one)
$ws.Rows("4:4").Insert()

$ws.Range("A4").Value = "Showroom Paused"
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "Project currently on hold pending strategic review"

# Insert a new row at row 7 for "Warehouse Paused" (after the now-shifted "Warehouse Target Date" row 6)
$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = "Warehouse Paused"
$ws.Range("B7").Value = "No"
$ws.Range("C7").Value = ""
